# Update "想去人数" (number of interested attendees) figures in both the
# "展览" sheet and the aggregated "全部类型" sheet to match the latest
# scrape output (gh-pages generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 657
$wsExhibit.Range("F12").Value = 10290
$wsExhibit.Range("F17").Value = 11862
$wsExhibit.Range("F18").Value = 12252

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 657
$wsAll.Range("F13").Value = 10290
$wsAll.Range("F18").Value = 11862
$wsAll.Range("F19").Value = 12252
